$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Users" column: header C1, first data value C2 = "Trade"
$ws.Range("C1").Value = "Users"
$ws.Range("C2").Value = "Trade"

# C1 gets the same bold+bordered header formatting as A1/B1.
$ws.Range("C1").Borders.Color = 0
$ws.Range("C1").Font.Bold = $true

# C2 gets a plain bordered style (no bold).
$ws.Range("C2").Borders.Color = 0

# Data validation: dropdown list of allowed user types on C2.
$ws.Range("C2").Validation.Add(3, 1, 1, '"Kam, Mmd, Trade"')

# Move the active selection to C2, matching the author's final cursor position.
$null = $ws.Range("C2").Select()
